$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1300.4
$ws.Range("I19").Value = 1270.375
$ws.Range("K19").Value = 1270.375
$ws.Range("M19").Value = -1095.375
$ws.Range("H51").Value = 7750.5
$ws.Range("J51").Value = 7750.5
$ws.Range("L51").Value = 7750.5
$ws.Range("N51").Value = -8718.5
$ws.Range("H99").Value = 2015.85
$ws.Range("I99").Value = 2393.5334
$ws.Range("K99").Value = 7180.600199999999
$ws.Range("M99").Value = -5682.600199999999
$ws.Range("H111").Value = 1007.1818
$ws.Range("J111").Value = 1081.3334
$ws.Range("L111").Value = 3244.0002
$ws.Range("N111").Value = -9378.0002
$ws.Range("H113").Value = 6953
$ws.Range("I113").Value = 5299.6665
$ws.Range("J113").Value = 7334.5386
$ws.Range("K113").Value = 5299.6665
$ws.Range("L113").Value = 7334.5386
$ws.Range("M113").Value = -2045.6665
$ws.Range("N113").Value = -13842.5386
$ws.Range("H129").Value = 1041.8572
$ws.Range("I129").Value = 715.5
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 2146.5
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = 2853.5
$ws.Range("N129").Value = -19000
$ws.Range("H137").Value = 38464724
$ws.Range("I137").Value = 66669396
$ws.Range("K137").Value = 200008188
$ws.Range("M137").Value = -200005638
$ws.Range("H138").Value = 4161.2607
$ws.Range("I138").Value = 1846.5
$ws.Range("J138").Value = 4508.475
$ws.Range("K138").Value = 5539.5
$ws.Range("L138").Value = 13525.425
$ws.Range("M138").Value = -399.5
$ws.Range("N138").Value = -23805.425

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 437.46667
$ws.Range("I5").Value = 152.125
$ws.Range("J5").Value = 763.5714
$ws.Range("K5").Value = 152.125
$ws.Range("L5").Value = 763.5714
$ws.Range("M5").Value = -40.125
$ws.Range("N5").Value = -987.5714
$ws.Range("H32").Value = 5455.017
$ws.Range("I32").Value = 3579.0364
$ws.Range("J32").Value = 31249.75
$ws.Range("K32").Value = 3579.0364
$ws.Range("L32").Value = 31249.75
$ws.Range("M32").Value = -3292.0364
$ws.Range("N32").Value = -31823.75
$ws.Range("H45").Value = 2391
$ws.Range("I45").Value = 2391
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2391
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2014
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 3784.0286
$ws.Range("I61").Value = 3047.5667
$ws.Range("J61").Value = 8202.799999999999
$ws.Range("K61").Value = 3047.5667
$ws.Range("L61").Value = 8202.799999999999
$ws.Range("M61").Value = -2835.5667
$ws.Range("N61").Value = -8626.799999999999
$ws.Range("H74").Value = 2094.0688
$ws.Range("I74").Value = 1489.7858
$ws.Range("K74").Value = 1489.7858
$ws.Range("M74").Value = -615.7858000000001
$ws.Range("H77").Value = 2094.0688
$ws.Range("I77").Value = 1489.7858
$ws.Range("K77").Value = 7448.929
$ws.Range("M77").Value = -3080.929
$ws.Range("H97").Value = 1464.8462
$ws.Range("I97").Value = 1093.125
$ws.Range("J97").Value = 2059.6
$ws.Range("K97").Value = 1093.125
$ws.Range("L97").Value = 2059.6
$ws.Range("M97").Value = -597.125
$ws.Range("N97").Value = -3051.6
$ws.Range("H114").Value = 60398.8
$ws.Range("J114").Value = 60398.8
$ws.Range("L114").Value = 60398.8
$ws.Range("N114").Value = -69076.8
$ws.Range("H131").Value = 64371.332
$ws.Range("J131").Value = 64371.332
$ws.Range("L131").Value = 64371.332
$ws.Range("N131").Value = -74451.33199999999
$ws.Range("H132").Value = 6296.2
$ws.Range("I132").Value = 3583
$ws.Range("J132").Value = 7459
$ws.Range("K132").Value = 10749
$ws.Range("L132").Value = 22377
$ws.Range("M132").Value = -8219
$ws.Range("N132").Value = -27437
$ws.Range("H136").Value = 3784.0286
$ws.Range("I136").Value = 3047.5667
$ws.Range("J136").Value = 8202.799999999999
$ws.Range("K136").Value = 9142.7001
$ws.Range("L136").Value = 24608.4
$ws.Range("M136").Value = -6592.7001
$ws.Range("N136").Value = -29708.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 437.46667
$ws.Range("I4").Value = 152.125
$ws.Range("J4").Value = 763.5714
$ws.Range("K4").Value = 152.125
$ws.Range("L4").Value = 763.5714
$ws.Range("M4").Value = -37.125
$ws.Range("N4").Value = -993.5714
$ws.Range("H20").Value = 4978.4165
$ws.Range("I20").Value = 4339.45
$ws.Range("K20").Value = 4339.45
$ws.Range("M20").Value = -4092.45
$ws.Range("H94").Value = 1722.5555
$ws.Range("I94").Value = 2079.1428
$ws.Range("J94").Value = 474.5
$ws.Range("K94").Value = 2079.1428
$ws.Range("L94").Value = 474.5
$ws.Range("M94").Value = -1628.1428
$ws.Range("N94").Value = -1376.5
$ws.Range("H99").Value = 1758.5834
$ws.Range("I99").Value = 1554.8182
$ws.Range("K99").Value = 1554.8182
$ws.Range("M99").Value = -56.81819999999993
$ws.Range("H105").Value = 22319.062
$ws.Range("I105").Value = 30739.857
$ws.Range("J105").Value = 15769.556
$ws.Range("K105").Value = 30739.857
$ws.Range("L105").Value = 15769.556
$ws.Range("M105").Value = -28992.857
$ws.Range("N105").Value = -19263.556
$ws.Range("H134").Value = 2067.2927
$ws.Range("I134").Value = 1401.2424
$ws.Range("J134").Value = 4814.75
$ws.Range("K134").Value = 4203.7272
$ws.Range("L134").Value = 14444.25
$ws.Range("M134").Value = -1668.7272
$ws.Range("N134").Value = -19514.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32443.572
$ws.Range("I31").Value = 1818.6316
$ws.Range("K31").Value = 1818.6316
$ws.Range("M31").Value = -1523.6316
$ws.Range("H34").Value = 32443.572
$ws.Range("I34").Value = 1818.6316
$ws.Range("K34").Value = 1818.6316
$ws.Range("M34").Value = -1616.6316
$ws.Range("H58").Value = 3234.353
$ws.Range("J58").Value = 5893.5835
$ws.Range("L58").Value = 5893.5835
$ws.Range("N58").Value = -6299.5835
$ws.Range("H99").Value = 1969.04
$ws.Range("I99").Value = 1701.625
$ws.Range("J99").Value = 2444.4443
$ws.Range("K99").Value = 1701.625
$ws.Range("L99").Value = 2444.4443
$ws.Range("M99").Value = -203.625
$ws.Range("N99").Value = -5440.4443
$ws.Range("H126").Value = 1969.04
$ws.Range("I126").Value = 1701.625
$ws.Range("J126").Value = 2444.4443
$ws.Range("K126").Value = 5104.875
$ws.Range("L126").Value = 7333.3329
$ws.Range("M126").Value = -2634.875
$ws.Range("N126").Value = -12273.3329
$ws.Range("H132").Value = 4029.0322
$ws.Range("I132").Value = 3489.6365
$ws.Range("K132").Value = 10468.9095
$ws.Range("M132").Value = -7938.9095
$ws.Range("H136").Value = 3234.353
$ws.Range("J136").Value = 5893.5835
$ws.Range("L136").Value = 17680.7505
$ws.Range("N136").Value = -22780.7505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 142859620
$ws.Range("J113").Value = 250002080
$ws.Range("L113").Value = 750006240
$ws.Range("N113").Value = -750010580
$ws.Range("H132").Value = 3671.8
$ws.Range("I132").Value = 3350
$ws.Range("K132").Value = 30150
$ws.Range("M132").Value = -27620
$ws.Range("H136").Value = 17546534
$ws.Range("I136").Value = 22223544
$ws.Range("J136").Value = 7747.5
$ws.Range("K136").Value = 66670632
$ws.Range("L136").Value = 23242.5
$ws.Range("M136").Value = -66665532
$ws.Range("N136").Value = -33442.5
$ws.Range("H139").Value = 3862.524
$ws.Range("I139").Value = 2549.9092
$ws.Range("J139").Value = 5306.4
$ws.Range("K139").Value = 7649.7276
$ws.Range("L139").Value = 15919.2
$ws.Range("M139").Value = -2509.7276
$ws.Range("N139").Value = -26199.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 507.05
$ws.Range("J2").Value = 1190.25
$ws.Range("L2").Value = 1190.25
$ws.Range("N2").Value = -1416.25
$ws.Range("H19").Value = 4999
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 4999
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 4999
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -5575
$ws.Range("H92").Value = 35140.57
$ws.Range("J92").Value = 35140.57
$ws.Range("L92").Value = 35140.57
$ws.Range("N92").Value = -38884.57
$ws.Range("H97").Value = 2216.4546
$ws.Range("I97").Value = 1631.4546
$ws.Range("J97").Value = 2801.4546
$ws.Range("K97").Value = 1631.4546
$ws.Range("L97").Value = 2801.4546
$ws.Range("M97").Value = -1135.4546
$ws.Range("N97").Value = -3793.4546
$ws.Range("H126").Value = 3647.2632
$ws.Range("I126").Value = 1944.9231
$ws.Range("K126").Value = 5834.7693
$ws.Range("M126").Value = -3364.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840
$ws.Range("H132").Value = 4386.0625
$ws.Range("I132").Value = 3590
$ws.Range("J132").Value = 7835.6665
$ws.Range("K132").Value = 10770
$ws.Range("L132").Value = 23506.9995
$ws.Range("M132").Value = -8240
$ws.Range("N132").Value = -28566.9995
$ws.Range("H136").Value = 3856.0566
$ws.Range("I136").Value = 2970.3928
$ws.Range("J136").Value = 4848
$ws.Range("K136").Value = 8911.178400000001
$ws.Range("L136").Value = 14544
$ws.Range("M136").Value = -6361.178400000001
$ws.Range("N136").Value = -19644

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 20997.5
$ws.Range("I30").Value = 12000
$ws.Range("K30").Value = 12000
$ws.Range("M30").Value = -11893
$ws.Range("H39").Value = 10016200
$ws.Range("J39").Value = 27000
$ws.Range("L39").Value = 27000
$ws.Range("N39").Value = -27826
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H126").Value = 3088
$ws.Range("I126").Value = 2402.2
$ws.Range("K126").Value = 7206.599999999999
$ws.Range("M126").Value = -4736.599999999999
$ws.Range("H132").Value = 3319.2
$ws.Range("I132").Value = 2630.682
$ws.Range("K132").Value = 7892.045999999999
$ws.Range("M132").Value = -5362.045999999999
